$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 51 (shifts rows 51:114 down to 52:115,
# carrying their existing formatting/values with them).
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new data point.
$ws.Range("A51").Value = 5
$ws.Range("B51").Value = "Macroferia Regional de Talca"
$ws.Range("C51").Value = "Maule"
$ws.Range("D51").Value = 45117
$ws.Range("E51").Value = 7
$ws.Range("F51").Value = 100112013
$ws.Range("G51").Value = "Alcachofa"
$ws.Range("H51").Value = "Madrigal"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 300
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = 15000
$ws.Range("N51").Value = "$/caja 40 unidades"
$ws.Range("O51").Value = "Provincia del Elquí"
$ws.Range("P51").Value = 375
$ws.Range("Q51").Value = 40
$ws.Range("R51").Value = "Hortaliza"

# Match the existing date-column style (numFmt "YYYY-MM-DD HH:MM:SS")
# used by every other row in column D.
$ws.Range("D51").NumberFormat = $ws.Range("D52").NumberFormat
